# Applies the "Add/update resource data for UWTranslationQuestions" edit:
#  1. Removes the "License Information" Heading2 paragraph.
#  2. Merges the license-text paragraph with the following
#     "This PDF version is provided under the same license." paragraph.
#  3. Rewrites the merged paragraph's body text/hyperlinks with the new
#     resource-license copy.

$d = $word.ActiveDocument

function Find-InRange($range, $needle) {
    $r = $range.Duplicate
    $null = $r.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $r
}

# A Range.Text assignment whose Start sits exactly on a run boundary picks
# up the formatting of the *preceding* character, which would make newly
# inserted plain text inherit e.g. Bold from an adjacent bold run. To avoid
# that we keep a one-character "buffer" from the untouched plain run,
# perform the bulk replacement strictly inside it, then delete the
# now-redundant buffer character.
function Set-TailText($startPos, $endPos, $newText) {
    $bulk = $d.Range($startPos + 1, $endPos)
    $bulk.Text = $newText
    $buffer = $d.Range($startPos, $startPos + 1)
    $buffer.Delete()
}

# ---------------------------------------------------------------------
# Step 1: delete the whole "License Information" heading paragraph
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "License Information`r") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# Step 2: locate the paragraph that holds the bold resource title and the
# license text, then merge it with the following paragraph
# ("This PDF version...") by deleting the paragraph mark between them.
# ---------------------------------------------------------------------
$licPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("is based on")) {
        $licPara = $p
        break
    }
}

$endOfLicPara = $licPara.Range.End
$d.Range($endOfLicPara - 1, $endOfLicPara).Delete()

# ---------------------------------------------------------------------
# Step 3: replace the bold run's text.
# ---------------------------------------------------------------------
$rBold = Find-InRange $licPara.Range "अनुवाद प्रश्न (unfoldingWord)"
$rBold.Text = "unfoldingWord® Translation Questions"

# ---------------------------------------------------------------------
# Step 4: replace everything from right after the bold run through the
# end of the paragraph with the new license copy (plain text for now --
# run boundaries are introduced afterwards).
# ---------------------------------------------------------------------
$boldEnd = $rBold.End
$paraEnd = $licPara.Range.End - 1   # stop short of the paragraph mark

$newTail = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. unfoldingWord® Translation Questions has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from unfoldingWord® Translation Questions © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual"

Set-TailText $boldEnd $paraEnd $newTail

Write-Output "done"
